# Bug Fixing Fitur Import Excel
# Insert a new header row at the top of the schedule export sheet so the
# exported columns line up with the database field names used by the
# import feature (id_jadwal, id_mhs, tanggal, jam_masuk, jam_keluar, ruangan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data down one row to make room for the new header row.
$ws.Rows.Item(1).Insert()

# Write the new header row that matches the DB column names expected by the
# Excel import feature.
$ws.Range("A1").Value = "id_jadwal"
$ws.Range("B1").Value = "id_mhs"
$ws.Range("C1").Value = "tanggal"
$ws.Range("D1").Value = "jam_masuk"
$ws.Range("E1").Value = "jam_keluar"
$ws.Range("F1").Value = "ruangan"

# Re-apply an explicit, import-friendly date format to the tanggal column
# now that it has shifted down to rows 2:16.
$ws.Range("C2:C16").NumberFormat = "yyyy\-mm\-dd"

# Keep the sheet set to portrait like the rest of the workbook's printouts.
$ws.PageSetup.Orientation = 1

# Restore the last active selection used while editing the sheet.
$ws.Range("G13").Select() | Out-Null
